$d = $word.ActiveDocument

# Table 2 (index 2, 1-based) holds the "Escuela" info block:
#   row 2: "Nombre completo" | <empty>
#   row 3: "Plan de Estudio"  | <empty>
#   row 4: "Año de ingreso"   | <empty>
$t = $d.Tables(2)

$t.Cell(2, 2).Range.Text = "Tomas Paz de la vega"
$t.Cell(3, 2).Range.Text = "Ing informatica"
$t.Cell(4, 2).Range.Text = "2022"
